# LISTA.xlsx edit script
# - strips the stray "Coordinador:" prefix from the two coordinator labels
#   (shared strings used throughout column E of Hoja1)
# - adds a new (currently empty) column F on Hoja1 with a custom width,
#   and leaves it selected with the view scrolled down one row, mirroring
#   the author's last interactive state before saving

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Fix the two mislabeled "Coordinador:..." entries in column E.
#    These are shared strings reused by many rows, so every cell that
#    used the old text has to be rewritten with the corrected text.
# ---------------------------------------------------------------------
$rowsZona44 = @(2,3,4,6,11,12,13,14,15,17,19,21,22,24,26,27,28,30,35,37)
$rowsZona45 = @(5,7,8,9,10,16,18,20,23,25,29,31,32,33,34,36,38,39)

foreach ($r in $rowsZona44) {
    $ws1.Range("E$r").Value = "COORDINADOR ZONA 4.4"
}

foreach ($r in $rowsZona45) {
    $ws1.Range("E$r").Value = "COORDINADOR ZONA 4.5"
}

# ---------------------------------------------------------------------
# 2) Introduce column F with its own width.
# ---------------------------------------------------------------------
$ws1.Columns.Item(6).ColumnWidth = 24.5

# ---------------------------------------------------------------------
# 3) Leave the sheet active with column F selected and the view
#    scrolled down a row, matching the saved selection state.
# ---------------------------------------------------------------------
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$ws1.Range("F2:F1048576").Select()
